$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 28, shifting existing rows 28-61 down by one.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new weekly data entry.
$ws.Cells.Item(28, 1).Value = 8
$ws.Cells.Item(28, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(28, 3).Value = "Coquimbo"
$ws.Cells.Item(28, 4).Value = 44586
$ws.Cells.Item(28, 5).Value = 4
$ws.Cells.Item(28, 6).Value = 100112030
$ws.Cells.Item(28, 7).Value = "Poroto granado"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 600
$ws.Cells.Item(28, 11).Value = 29000
$ws.Cells.Item(28, 12).Value = 30000
$ws.Cells.Item(28, 13).Value = 29500
$ws.Cells.Item(28, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(28, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(28, 16).Value = 1180
$ws.Cells.Item(28, 17).Value = 25
$ws.Cells.Item(28, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same style/number format as other date cells in column D.
$ws.Cells.Item(28, 4).NumberFormat = $ws.Cells.Item(29, 4).NumberFormat
